$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Product Id"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Current Quantity"
$ws.Range("D1").Value = "Latest Supplier Name"
$ws.Range("E1").Value = "Latest Supplier Unit Price"
$ws.Range("F1").Value = "Latest Supplier Date"
$ws.Range("G1").Value = "Minimum Supplier Name"
$ws.Range("H1").Value = "Minimum Supplier Unit Price"
$ws.Range("I1").Value = "Minimum Supplier Date"
